$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.600.23'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.808.60'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '225.65'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.598'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.76%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '37.38'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +6.83%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.292'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.95%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0682'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.070.80'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.32'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.808.99'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '34.518.53'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.75'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '243.57'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0777'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.88%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.23'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.13'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.06%  '
$ws.Range('E24').Value = '  +4.80%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '171.79'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.88'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.32'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.75%  '
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.94'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.31%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.82'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0518'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.66%  '
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.365.97'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.17%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.655'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.12%  '
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('E38').Value = '  -5.04%  '
$ws.Range('E39').Value = '  -1.49%  '
$ws.Range('E40').Value = '  +1.63%  '
$ws.Range('E41').Value = '  -1.46%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '81.00'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.17%  '
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('E44').Value = '  +5.57%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.83'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.970.83'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.82'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '102.88'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.01%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0₆0121'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -7.33%  '
